$d = $word.ActiveDocument

# 1. Merge "to " + _GoBack bookmark + "manually set" into a single run "to manually set".
#    Find/Replace spans the (now stranded) _GoBack bookmark and collapses it away, matching
#    the diff's removal of bookmarkStart/End id=6 at this spot.
$d.Content.Find.Execute("to manually set", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "to manually set", 2)

# 2. Append the dwell-time-limit sentence to the end of the Arbitrary paragraph.
$d.Content.Find.Execute("tune to the wavelength, or latency in the USB communications.", `
                         $false, $false, $false, $false, $false, $true, 1, $false, `
                         "tune to the wavelength, or latency in the USB communications. Dwell time is limited to between 0 and 10000 milliseconds (10 seconds).", `
                         2)

# 3. Append the dwell-time-limit sentence to the end of the Ordered paragraph, and add a new
#    empty paragraph right after it (before the pre-existing trailing empty paragraph).
$d.Content.Find.Execute("by setting Start higher than Stop and using a negative Step. ", `
                         $false, $false, $false, $false, $false, $true, 1, $false, `
                         "by setting Start higher than Stop and using a negative Step. Dwell time is limited to between 0 and 10000 milliseconds (10 seconds).^p", `
                         2)

# 4. Word keeps a "_GoBack" bookmark at the site of the most recent edit. Re-create it inside
#    the Arbitrary paragraph, right before "ecify" in "specify" (matching the author's last
#    edit location after typing the new sentences above).
$r = $d.Content
$r.Find.Execute("allows the user to sp", $false, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0)
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)
